# Update the "Export" sheet's account/name/balance rows (rows 2-11) to
# reflect the new Saldo export snapshot. Net effect vs. the previous
# snapshot:
#   - 4 new rows added (EDUARDO, VENIA, and refreshed RICARDO / THIAGO
#     balances) ahead of the existing block
#   - MARCELO's balance updated from 128835.58 to 2000
#   - CRISTIANO added as a new row
#   - CRISTINA, JOSE and RAFAEL rows removed (no longer present)
#   - RODRIGO / GUILHERME / CARLOS / GISELA rows shift down but keep
#     their values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("004461070", "EDUARDO",   145045.2),
    @("004813166", "VENIA",      52636.3),
    @("005046805", "RICARDO",    51658.83),
    @("005064129", "THIAGO",     20349.02),
    @("003641655", "MARCELO",    2000),
    @("004940560", "CRISTIANO",  1113.24),
    @("004392159", "RODRIGO",    900.21),
    @("004574428", "GUILHERME",  745.08),
    @("004488571", "CARLOS",     440.36),
    @("004322719", "GISELA",     276.97)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $acct = $data[$i][0]
    $name = $data[$i][1]
    $bal = $data[$i][2]

    $acctCell = $ws.Cells.Item($row, 1)
    $nameCell = $ws.Cells.Item($row, 2)
    $balCell = $ws.Cells.Item($row, 3)

    # Only touch the account-number cell if it doesn't already hold the
    # right text. Leading with an apostrophe keeps the numeric-looking
    # string (e.g. "004461070") as text with its leading zeros intact
    # instead of Excel silently coercing it to a number; skipping the
    # write when unchanged also avoids flipping an untouched cell's
    # style (quote-prefix) for no reason.
    if ($acctCell.Text -ne $acct) {
        $acctCell.Value = "'" + $acct
    }
    if ($nameCell.Text -ne $name) {
        $nameCell.Value = $name
    }
    if ($balCell.Value2 -ne $bal) {
        $balCell.Value = $bal
    }
}
